# Update the "Prix Spot" sheet: add a new day column (R) for 01-jul with
# the same header styling as the existing date columns, then fill in the
# hourly price values for that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the last existing header cell (Q1) onto the new
# header cell (R1) so it keeps the same bold/centered/bordered style,
# then set its text.
$ws.Range("Q1").Copy()
$ws.Range("R1").PasteSpecial(-4122)
$ws.Range("R1").Value = "01-jul"

# Hourly prices for 01-jul, row 2 (00-01) through row 25 (23-24).
$values = @(
    111.28,
    95.41,
    89.09999999999999,
    88.08,
    90.01000000000001,
    96.64,
    114.97,
    120.06,
    114.74,
    93.19,
    94.91,
    80.7,
    75.11,
    65.53,
    68.95,
    78.16,
    90.8,
    102.65,
    150,
    223.54,
    235,
    215,
    199.3,
    129.78
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 18).Value = $values[$i]
}
